$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 header values changed
$ws.Range("B1").Value = 16
$ws.Range("C1").Value = 20
$ws.Range("D1").Value = 16
$ws.Range("E1").Value = 20

# Row 2: D2 removed, B2 and C2 added
$ws.Range("D2").Value = $null
$ws.Range("B2").Value = 30.373805491377226
$ws.Range("C2").Value = 27.783653487400766

# Row 3: B3 removed, C3 changed
$ws.Range("B3").Value = $null
$ws.Range("C3").Value = 25.230453533488497

# Selection changed to B1:E3
$ws.Range("B1:E3").Select()
